# Updates the cryptos list figures (Price / Volume(1h)) to reflect the
# latest scrape, as produced by the "Updated cryptos list ... with GitHub
# Actions" workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new value.
# Only cells that actually changed are listed here.
$changes = @{
    2  = @{ D = "68.123.42";  E = "  +0.61%  " }
    3  = @{ D = "3.792.43";   E = "  -0.22%  " }
    4  = @{ D = "1.00";       E = "  +0.01%  " }
    5  = @{ D = "600.78";     E = "  +0.73%  " }
    6  = @{ D = "165.44";     E = "  -1.05%  " }
    7  = @{ E = "  +0.06%  " }
    8  = @{ E = "  -0.73%  " }
    9  = @{ E = "  -1.27%  " }
    10 = @{ D = "0.451";      E = "  +0.36%  " }
    11 = @{ E = "  +3.12%  " }
    12 = @{ D = "0.0000249";  E = "  -1.45%  " }
    13 = @{ D = "35.76";      E = "  -0.91%  " }
    14 = @{ D = "4.430.00";   E = "  -0.21%  " }
    15 = @{ D = "3.821.33";   E = "  +0.48%  " }
    16 = @{ D = "68.139.66";  E = "  +0.65%  " }
    17 = @{ D = "18.45";      E = "  -0.70%  " }
    18 = @{ E = "  +2.17%  " }
    19 = @{ D = "7.07";       E = "  -0.56%  " }
    20 = @{ D = "461.19";     E = "  -0.02%  " }
    21 = @{ D = "9.70";       E = "  -2.26%  " }
    22 = @{ D = "0.699";      E = "  -0.44%  " }
    23 = @{ D = "0.0000149";  E = "  -3.30%  " }
    24 = @{ D = "82.89";      E = "  -0.76%  " }
    25 = @{ D = "12.06";      E = "  -0.28%  " }
    26 = @{ E = "  +0.11%  " }
    27 = @{ E = "  -0.72%  " }
    28 = @{ D = "9.98";       E = "  -0.26%  " }
    29 = @{ D = "3.942.34";   E = "  -0.19%  " }
    30 = @{ D = "7.37";       E = "  +1.76%  " }
    31 = @{ D = "2.64";       E = "  -5.26%  " }
    32 = @{ D = "2.22";       E = "  -1.88%  " }
    33 = @{ D = "29.29";      E = "  -1.43%  " }
    34 = @{ E = "  +0.03%  " }
    35 = @{ D = "8.99";       E = "  -0.96%  " }
    36 = @{ D = "0.1000";     E = "  -0.04%  " }
    37 = @{ E = "  +0.45%  " }
    38 = @{ D = "3.29";       E = "  -2.86%  " }
    39 = @{ E = "  -0.06%  " }
    40 = @{ D = "0.985";      E = "  -0.94%  " }
    41 = @{ D = "1.00" }
    42 = @{ E = "  -0.01%  " }
    # Arweave and TheGraph swap ranking positions, each keeping its own
    # updated price / volume figures.
    43 = @{ B = "Arweave";  C = "https://coinranking.com/coin/7XWg41D1+arweave-ar";      D = "43.63"; E = "  -0.52%  " }
    44 = @{ B = "TheGraph"; C = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt";    D = "0.300"; E = "  +0.48%  " }
    45 = @{ D = "47.36";      E = "  -1.52%  " }
    46 = @{ D = "151.66";     E = "  +1.06%  " }
    47 = @{ D = "8.35";       E = "  +0.24%  " }
    48 = @{ E = "  +2.77%  " }
    49 = @{ E = "  +6.77%  " }
    50 = @{ D = "393.53";     E = "  -0.09%  " }
    51 = @{ D = "26.67";      E = "  +1.03%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        # Force the value to be stored as text (these cells hold things
        # like prices "1.00" / "600.78" which Excel would otherwise
        # helpfully "convert" into real numbers), then strip the
        # formatting changes that setting NumberFormat introduced so the
        # cell's style matches the rest of the untouched column.
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        $cell.ClearFormats()
    }
}
